$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title / header shared text (volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/9/2024  Through  9/15/2024"

# --- Fix number formats for cells changing between text placeholder ("0"/"***.*") and real numbers ---
# Reference cells that keep a stable style throughout this edit, used as format-paste sources.
$fmtText = $ws.Range("C23")   # style 14: text / General (used for "0" and "***.*" placeholders)
$fmtCount = $ws.Range("I15")  # style 16: #,##0 integer count format
$fmtPct = $ws.Range("M15")    # style 15: #,##0.0;"-"#,##0.0 percent-change format

$fmtCount.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1

$fmtPct.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100

$fmtCount.Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 2

$fmtText.Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "0"

$fmtCount.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 2

$fmtPct.Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100

$fmtCount.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1

$fmtPct.Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100

$fmtText.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = "0"

$fmtText.Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"

$fmtCount.Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1

$fmtPct.Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100

$fmtCount.Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1

$fmtPct.Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100

$fmtCount.Copy()
$ws.Range("J29").PasteSpecial(-4122)
$ws.Range("J29").Value = 1

$fmtPct.Copy()
$ws.Range("K29").PasteSpecial(-4122)
$ws.Range("K29").Value = 0

$fmtCount.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1

$fmtPct.Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100

$fmtCount.Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1

$fmtPct.Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100

$fmtCount.Copy()
$ws.Range("J30").PasteSpecial(-4122)
$ws.Range("J30").Value = 1

$fmtPct.Copy()
$ws.Range("K30").PasteSpecial(-4122)
$ws.Range("K30").Value = 0

# --- Plain value updates (style unchanged) ---
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 40
$ws.Range("L15").Value = -36.363636363636
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = -70.588235294117
$ws.Range("I16").Value = 70
$ws.Range("J16").Value = 88
$ws.Range("K16").Value = -20.454545454545
$ws.Range("L16").Value = -27.835051546391
$ws.Range("M16").Value = -22.222222222222
$ws.Range("N16").Value = -83.529411764705
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 6
$ws.Range("H17").Value = -14.285714285714
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 11.25
$ws.Range("L17").Value = -23.931623931623
$ws.Range("N17").Value = -52.150537634408
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 42.857142857142
$ws.Range("I18").Value = 102
$ws.Range("J18").Value = 131
$ws.Range("K18").Value = -22.137404580152
$ws.Range("L18").Value = -44.262295081967
$ws.Range("M18").Value = -43.956043956044
$ws.Range("N18").Value = -87.666263603385
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 2.083333333333
$ws.Range("I19").Value = 452
$ws.Range("J19").Value = 476
$ws.Range("K19").Value = -5.042016806722
$ws.Range("L19").Value = 7.619047619047
$ws.Range("M19").Value = 121.56862745098
$ws.Range("N19").Value = 76.5625
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 75
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = -37.5
$ws.Range("L20").Value = -41.860465116279
$ws.Range("M20").Value = -28.571428571428
$ws.Range("N20").Value = -88.408037094281
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 15
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = -11.702127659574
$ws.Range("I21").Value = 795
$ws.Range("J21").Value = 900
$ws.Range("K21").Value = -11.666666666666
$ws.Range("L21").Value = -17.014613778705
$ws.Range("M21").Value = 24.21875
$ws.Range("N21").Value = -66.213344666383
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 21
$ws.Range("K23").Value = -4.761904761904
$ws.Range("L23").Value = 5.263157894736
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 11.764705882352
$ws.Range("F24").Value = 74
$ws.Range("H24").Value = 5.714285714285
$ws.Range("I24").Value = 701
$ws.Range("J24").Value = 651
$ws.Range("K24").Value = 7.680491551459
$ws.Range("L24").Value = -2.231520223152
$ws.Range("M24").Value = 70.559610705596
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 62.5
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 430
$ws.Range("J25").Value = 349
$ws.Range("K25").Value = 23.209169054441
$ws.Range("L25").Value = 5.651105651105
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 41.176470588235
$ws.Range("I26").Value = 194
$ws.Range("J26").Value = 170
$ws.Range("K26").Value = 14.117647058823
$ws.Range("L26").Value = -4.433497536945
$ws.Range("M26").Value = 24.358974358974
$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 14.285714285714
$ws.Range("L27").Value = -38.461538461538
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 27
$ws.Range("K28").Value = -15.625
$ws.Range("L28").Value = 17.391304347826

$excel.CutCopyMode = $false
